# "updated Incidence to remove single pollutant lines with inc = 0"
#
# The Incidence sheet lists, per Source ID / Pollutant combination, the
# incidence rate of a given emission. Rows whose Source ID is "Total" are
# aggregate (all-sources) rows; every other Source ID (e.g. "CEOT0001")
# describes a single emission source / pollutant. This edit drops the
# "single pollutant" rows (Source ID <> "Total") whose Incidence (column D)
# is 0, while leaving the aggregate "Total" rows untouched even when their
# Incidence is also 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last populated row in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Walk the data rows bottom-to-top so that deleting a row doesn't shift the
# row numbers of the ones still left to examine.
for ($r = $lastRow; $r -ge 2; $r--) {
    $sourceId = $ws.Cells.Item($r, 1).Text
    $incidence = $ws.Cells.Item($r, 4).Text
    if ($sourceId -ne "Total" -and $incidence -eq "0") {
        $ws.Rows.Item($r).Delete()
    }
}
